# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 273
$ws1.Range("F4").Value = 1016
$ws1.Range("F5").Value = 557

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 273
$ws4.Range("F4").Value = 1016
$ws4.Range("F6").Value = 557
